# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the "Casos activos" (active cases) values between the Huesca/Huelva
# rows: row 53 (previously Huelva=72) becomes Huesca=0, row 54 (previously
# Huesca=0) becomes Huelva=72.
$ws.Range("A53").Value = "Huesca"
$ws.Range("C53").Value = 0

$ws.Range("A54").Value = "Huelva"
$ws.Range("C54").Value = 72

# Update the "last updated" timestamp string.
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 05:46"
